$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (pushes existing rows 13+ down by one)
$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C13").Value = "Arica y Parinacota"
$ws.Range("D13").Value = 45204
$ws.Range("E13").Value = 15
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100101
$ws.Range("H13").Value = "Berries"
$ws.Range("I13").Value = 100101007
$ws.Range("J13").Value = "Kiwi"
$ws.Range("K13").Value = "Hayward"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 300
$ws.Range("N13").Value = 23000
$ws.Range("O13").Value = 24000
$ws.Range("P13").Value = 23500
$ws.Range("Q13").Value = "`$/bandeja 10 kilos"
$ws.Range("R13").Value = "Región de O'Higgins"
$ws.Range("S13").Value = 2350
$ws.Range("T13").Value = 10

# Insert a second new row at position 24 (after the first insert has shifted
# the old row 22 down to row 23; pushes rows 24+ down by one more)
$ws.Rows.Item(24).Insert()

$ws.Range("A24").Value = 1
$ws.Range("B24").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C24").Value = "Arica y Parinacota"
$ws.Range("D24").Value = 45223
$ws.Range("E24").Value = 15
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100101
$ws.Range("H24").Value = "Berries"
$ws.Range("I24").Value = 100101007
$ws.Range("J24").Value = "Kiwi"
$ws.Range("K24").Value = "Hayward"
$ws.Range("L24").Value = "Especial"
$ws.Range("M24").Value = 250
$ws.Range("N24").Value = 24000
$ws.Range("O24").Value = 25000
$ws.Range("P24").Value = 24600
$ws.Range("Q24").Value = "`$/bandeja 10 kilos"
$ws.Range("R24").Value = "Región de O'Higgins"
$ws.Range("S24").Value = 2460
$ws.Range("T24").Value = 10
